# main_RF translated to jupyter notebook
# Reposition/resize the two pictures on slide 2 (content-placeholder image and the
# free Picture 3) to their new EMU coordinates.
#
# Shape.Left/Top/Width/Height round-trip through a single-precision (float32)
# point value, so the literals below are the points values (pre-computed) whose
# float32 representation reconstructs exactly the target EMU, instead of the
# naive "emu / (914400/72)" which can drift by a single EMU after the
# float64 -> float32 -> EMU trip.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- "Content Placeholder 3" picture: off (838200,1968090)/ext (4686300,3400425)
#     -> off (1262744,2078231)/ext (3088741,2241221) ---
$pic1 = $s.Shapes.Item("Content Placeholder 3")
$pic1.Left   = 99.42870330810547
$pic1.Top    = 163.64027404785156
$pic1.Width  = 243.20799255371094
$pic1.Height = 176.47413635253906

# --- "Picture 3": off (5843587,1968090)/ext (4619625,3295650)
#     -> off (4676638,2078231)/ext (3088742,2203515) ---
$pic2 = $s.Shapes.Item("Picture 3")
$pic2.Left   = 368.2392578125
$pic2.Top    = 163.64027404785156
$pic2.Width  = 243.20806884765625
$pic2.Height = 173.50515747070312
